# Workbook / active sheet references
$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet   # "Table_1"

# -----------------------------------------------------------------
# 1) Clean up the handful of cells on "Table_1" that only ever held
#    an empty inline string (no real content) - they should not be
#    present at all after the edit.
# -----------------------------------------------------------------
$ws1.Range("B2").ClearContents()
$ws1.Range("A3").ClearContents()
$ws1.Range("B26").ClearContents()
$ws1.Range("B37").ClearContents()

# -----------------------------------------------------------------
# 2) Add a brand new worksheet "Table_2" right after "Table_1" and
#    fill it with the capital-adequacy ratio table.
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "6.0%"
$ws2.Range("B2").Style = "Normal"
$ws2.Range("C2").Value = "ratAdequacy1stDegreeCap"
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "5.0%"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "11.0%"
$ws2.Range("B3").Style = "Normal"
$ws2.Range("C3").Value = "ratAdequacyCumulativeCap"
$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "9.0%"
$ws2.Range("D3").Style = "Normal"

$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "ratLeverage"
$ws2.Range("D4").Value = "minimum 4%"

# Re-use the existing bold/centered/bordered header style (the one
# applied to A1:B1 on "Table_1") for the header row of "Table_2"
# instead of inventing a brand-new style entry.
$ws1.Range("A1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
